$d = $word.ActiveDocument

# Locate the phrase "an optional grouping of " (part of the "Routines" bullet
# under Functional Requirements) so we can append a brand-new run
# ("Activities") right after it, leaving the existing runs untouched —
# matches the target OOXML: a new <w:r><w:t>Activities</w:t></w:r> appended
# after the existing " \u2013 an optional grouping of " run.
$searchText = "an optional grouping of "

$rng = $d.Content
$found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found -eq $true) {
    # $rng now spans exactly the found text; move to its end and insert
    # a new run containing "Activities" right after it.
    $rng.Collapse(0)  # wdCollapseEnd
    $rng.InsertAfter("Activities")
}
